$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.145393333333334
$ws.Range("H2").Value = 12.43618
$ws.Range("I2").Value = 0.1621900462138432
$ws.Range("J2").Value = 0.1621900462138432
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.750436666666666
$ws.Range("N2").Value = 14.25131
$ws.Range("O2").Value = 0.07745299862590357
$ws.Range("P2").Value = 0.07745299862590359
$ws.Range("Q2").Value = 19.69242848842222
$ws.Range("R2").Value = 177.2318563958
$ws.Range("S2").Value = 0.01256210542653604
$ws.Range("T2").Value = 0.01256210542653604

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.145393333333334
$ws.Range("H3").Value = 12.43618
$ws.Range("I3").Value = 0.1621900462138432
$ws.Range("J3").Value = 0.1621900462138432
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 34.05277366666667
$ws.Range("N3").Value = 102.158321
$ws.Range("O3").Value = 0.5552098927072401
$ws.Range("P3").Value = 0.5552098927072401
$ws.Range("Q3").Value = 141.1621409393089
$ws.Range("R3").Value = 1270.45926845378
$ws.Range("S3").Value = 0.09004951815657021
$ws.Range("T3").Value = 0.0900495181565702

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 4.145393333333334
$ws.Range("H4").Value = 12.43618
$ws.Range("I4").Value = 0.1621900462138432
$ws.Range("J4").Value = 0.1621900462138432
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.52994333333334
$ws.Range("N4").Value = 67.58983000000001
$ws.Range("O4").Value = 0.3673371086668564
$ws.Range("P4").Value = 0.3673371086668564
$ws.Range("Q4").Value = 93.39547689437779
$ws.Range("R4").Value = 840.5592920494001
$ws.Range("S4").Value = 0.05957842263073698
$ws.Range("T4").Value = 0.05957842263073698

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.772999666666666
$ws.Range("H5").Value = 20.318999
$ws.Range("I5").Value = 0.2649961151116367
$ws.Range("J5").Value = 0.2649961151116367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.750436666666666
$ws.Range("N5").Value = 14.25131
$ws.Range("O5").Value = 0.07745299862590357
$ws.Range("P5").Value = 0.07745299862590359
$ws.Range("Q5").Value = 32.17470595985444
$ws.Range("R5").Value = 289.5723536386899
$ws.Range("S5").Value = 0.02052474373961138
$ws.Range("T5").Value = 0.02052474373961138

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.772999666666666
$ws.Range("H6").Value = 20.318999
$ws.Range("I6").Value = 0.2649961151116367
$ws.Range("J6").Value = 0.2649961151116367
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 34.05277366666667
$ws.Range("N6").Value = 102.158321
$ws.Range("O6").Value = 0.5552098927072401
$ws.Range("P6").Value = 0.5552098927072401
$ws.Range("Q6").Value = 230.6394246934088
$ws.Range("R6").Value = 2075.754822240679
$ws.Range("S6").Value = 0.1471284646389673
$ws.Range("T6").Value = 0.1471284646389673

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.772999666666666
$ws.Range("H7").Value = 20.318999
$ws.Range("I7").Value = 0.2649961151116367
$ws.Range("J7").Value = 0.2649961151116367
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.52994333333334
$ws.Range("N7").Value = 67.58983000000001
$ws.Range("O7").Value = 0.3673371086668564
$ws.Range("P7").Value = 0.3673371086668564
$ws.Range("Q7").Value = 152.5952986866856
$ws.Range("R7").Value = 1373.35768818017
$ws.Range("S7").Value = 0.09734290673305807
$ws.Range("T7").Value = 0.09734290673305808

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Bmp6"
$ws.Range("C8").Value = "Bmpr1a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.640471
$ws.Range("H8").Value = 43.921413
$ws.Range("I8").Value = 0.5728138386745201
$ws.Range("J8").Value = 0.5728138386745202
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.750436666666666
$ws.Range("N8").Value = 14.25131
$ws.Range("O8").Value = 0.07745299862590357
$ws.Range("P8").Value = 0.07745299862590359
$ws.Range("Q8").Value = 69.54863025566999
$ws.Range("R8").Value = 625.93767230103
$ws.Range("S8").Value = 0.04436614945975616
$ws.Range("T8").Value = 0.04436614945975617

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Bmp6"
$ws.Range("C9").Value = "Bmpr1a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.640471
$ws.Range("H9").Value = 43.921413
$ws.Range("I9").Value = 0.5728138386745201
$ws.Range("J9").Value = 0.5728138386745202
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 34.05277366666667
$ws.Range("N9").Value = 102.158321
$ws.Range("O9").Value = 0.5552098927072401
$ws.Range("P9").Value = 0.5552098927072401
$ws.Range("Q9").Value = 498.548645336397
$ws.Range("R9").Value = 4486.937808027573
$ws.Range("S9").Value = 0.3180319099117026
$ws.Range("T9").Value = 0.3180319099117027

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bmp6"
$ws.Range("C10").Value = "Bmpr1a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.640471
$ws.Range("H10").Value = 43.921413
$ws.Range("I10").Value = 0.5728138386745201
$ws.Range("J10").Value = 0.5728138386745202
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.52994333333334
$ws.Range("N10").Value = 67.58983000000001
$ws.Range("O10").Value = 0.3673371086668564
$ws.Range("P10").Value = 0.3673371086668564
$ws.Range("Q10").Value = 329.84898200331
$ws.Range("R10").Value = 2968.64083802979
$ws.Range("S10").Value = 0.2104157793030613
$ws.Range("T10").Value = 0.2104157793030614
